$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# Rotate data among rows 2, 3, 5:
#   new row2 = old row5
#   new row3 = old row2
#   new row5 = old row3
# Columns A, B, D, E, F, G, H, Q, R change values.
# Column L presence moves from row5 -> row2 (emptied on row5).
# Column AC presence/value moves from row5 -> row2 (emptied on row5).

# Capture old values first (rows 2, 3, 5) for columns that differ.
$cols = @("A","B","D","E","F","G","H","Q","R")

$old2 = @{}
$old3 = @{}
$old5 = @{}
foreach ($c in $cols) {
    $old2[$c] = $ws.Range("$c`2").Value2
    $old3[$c] = $ws.Range("$c`3").Value2
    $old5[$c] = $ws.Range("$c`5").Value2
}
$oldAC5 = $ws.Range("AC5").Value2

# Apply rotation: row2 <- old5, row3 <- old2, row5 <- old3
foreach ($c in $cols) {
    $ws.Range("$c`2").Value = $old5[$c]
    $ws.Range("$c`3").Value = $old2[$c]
    $ws.Range("$c`5").Value = $old3[$c]
}

# Column L: row5 had an (empty) inline string cell, row2 did not.
# Move that empty-string presence to L2 (copy an existing empty text
# cell so the cell record actually materializes), and clear L5.
$ws.Range("I2").Copy($ws.Range("L2"))
$ws.Range("L5").ClearContents()

# Column AC: move "På murken låga" comment from row5 to row2.
$ws.Range("AC2").Value = $oldAC5
$ws.Range("AC5").Value = $null
